# #5: property boat&car done
# Sheet3 = 汽車 (car). Expand the header row and the single data row so the
# sheet carries the same trailing metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index) that
# the other property sheets (land, house, deposit, stock, ...) already have,
# and give the first two columns real header labels ("name", "capacity")
# instead of duplicating the data row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Row 1 (header labels) ----
# B1/C1 currently hold copies of the data values ("NISSAN"/3000) instead of
# header text - replace them with the proper column names.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

# New trailing header cells - copy style from an existing header cell first
# so they pick up the same (bold/bordered) formatting as B1:G1.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("G1").Copy($ws.Range("I1"))
$ws.Range("G1").Copy($ws.Range("J1"))
$ws.Range("G1").Copy($ws.Range("K1"))
$ws.Range("G1").Copy($ws.Range("L1"))
$ws.Range("G1").Copy($ws.Range("M1"))
$ws.Range("G1").Copy($ws.Range("N1"))

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ---- Row 2 (data) ----
# B2:G2 keep their existing values (name/capacity/owner/register_date/
# register_reason/acquire_value) - only the new trailing columns are added.
$ws.Range("G2").Copy($ws.Range("H2"))
$ws.Range("G2").Copy($ws.Range("I2"))
$ws.Range("G2").Copy($ws.Range("J2"))
$ws.Range("G2").Copy($ws.Range("K2"))
$ws.Range("G2").Copy($ws.Range("L2"))
$ws.Range("G2").Copy($ws.Range("M2"))
$ws.Range("G2").Copy($ws.Range("N2"))

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# Force J2 to stay plain text ("2012-04-18") instead of being auto-parsed
# into a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-18"

$ws.Range("K2").Value = "邱議瑩"
$ws.Range("L2").Value = 913
$ws.Range("M2").Value = "tmped121"
$ws.Range("N2").Value = 29
